$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.503.13"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.618.02"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.57"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.02"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0886"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.848.84"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").Value = "1.618.75"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.548"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.41"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "27.506.00"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.51"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "0.0₃0720"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.88"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +7.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.08"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").Value = "1.451.76"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.34"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.931"
$ws.Range("E37").Value = "  +5.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.561"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.04"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.48"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.40"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.22"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("D47").Value = "1.759.71"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.13"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0985"
$ws.Range("E51").Value = "  -0.34%  "
